# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition list) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1943   # was 1932
$ws1.Range("F4").Value = 837    # was 833
$ws1.Range("F5").Value = 964    # was 938
$ws1.Range("F6").Value = 336    # was 278

# --- Sheet "全部类型" (all types list) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1943   # was 1932
$ws4.Range("F5").Value = 837    # was 833
$ws4.Range("F6").Value = 964    # was 938
$ws4.Range("F7").Value = 336    # was 278
